$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -0.4429594313435204
$ws.Cells.Item(2, 2).Value = -0.3992645103509139
$ws.Cells.Item(3, 1).Value = -0.4103409201240907
$ws.Cells.Item(3, 2).Value = -0.3765695994791023
$ws.Cells.Item(4, 1).Value = -0.3587143775834631
$ws.Cells.Item(4, 2).Value = -0.4018854641427108
$ws.Cells.Item(5, 1).Value = -0.1735750763153188
$ws.Cells.Item(5, 2).Value = -0.2816664731996503
$ws.Cells.Item(6, 1).Value = 0.03489176074543628
$ws.Cells.Item(6, 2).Value = -0.0935226103747592
$ws.Cells.Item(7, 1).Value = -0.3664722548637366
$ws.Cells.Item(7, 2).Value = -0.3599127665629931
$ws.Cells.Item(8, 1).Value = -0.4589633803383874
$ws.Cells.Item(8, 2).Value = -0.4070203543070836
$ws.Cells.Item(9, 1).Value = -0.5089773767051985
$ws.Cells.Item(9, 2).Value = -0.4273772976278596
$ws.Cells.Item(10, 1).Value = -0.5498162369676912
$ws.Cells.Item(10, 2).Value = -0.374545635298684
$ws.Cells.Item(11, 1).Value = -0.2013791479754209
$ws.Cells.Item(11, 2).Value = -0.310461465834296
$ws.Cells.Item(12, 1).Value = -0.2787778993590344
$ws.Cells.Item(12, 2).Value = -0.3021546202850426
$ws.Cells.Item(13, 1).Value = 0.1403567773520446
$ws.Cells.Item(13, 2).Value = 0.009251756584211377
$ws.Cells.Item(14, 1).Value = -0.1012906188007725
$ws.Cells.Item(14, 2).Value = -0.1319372724767337
$ws.Cells.Item(15, 1).Value = -0.1312577762006944
$ws.Cells.Item(15, 2).Value = -0.05067498419515562
$ws.Cells.Item(16, 1).Value = -0.2223208213173879
$ws.Cells.Item(16, 2).Value = -0.1103074485298595
$ws.Cells.Item(17, 1).Value = -0.02169987064437204
$ws.Cells.Item(17, 2).Value = 0.07646202848093735
$ws.Cells.Item(18, 1).Value = 0.01688094221336267
$ws.Cells.Item(18, 2).Value = 0.0957865711481915
$ws.Cells.Item(19, 1).Value = 0.04712566691634587
$ws.Cells.Item(19, 2).Value = 0.06285428990819035
$ws.Cells.Item(20, 1).Value = -0.1099620914004087
$ws.Cells.Item(20, 2).Value = -0.02390915371342364
$ws.Cells.Item(21, 1).Value = -0.1034741663147644
$ws.Cells.Item(21, 2).Value = -0.07545065315564811
$ws.Cells.Item(22, 1).Value = 0.05953571996208409
$ws.Cells.Item(22, 2).Value = 0.0883658863476722
$ws.Cells.Item(23, 1).Value = 0.1464629830698457
$ws.Cells.Item(23, 2).Value = 0.0691661075748489
$ws.Cells.Item(24, 1).Value = 0.4928127554208245
$ws.Cells.Item(24, 2).Value = 0.352530925714285
$ws.Cells.Item(25, 1).Value = 0.1871713219774291
$ws.Cells.Item(25, 2).Value = 0.1530511597418706
$ws.Cells.Item(26, 1).Value = 0.1983283978110824
$ws.Cells.Item(26, 2).Value = 0.155439310795237
$ws.Cells.Item(27, 1).Value = 0.1651000927412564
$ws.Cells.Item(27, 2).Value = 0.1400233657687214
$ws.Cells.Item(28, 1).Value = 0.2526736426165058
$ws.Cells.Item(28, 2).Value = 0.1833501533540975
$ws.Cells.Item(29, 1).Value = 0.4968485983507317
$ws.Cells.Item(29, 2).Value = 0.4141418490653193
$ws.Cells.Item(30, 1).Value = 0.2010148268514706
$ws.Cells.Item(30, 2).Value = 0.1824566527411436
$ws.Cells.Item(31, 1).Value = 0.1452498844205191
$ws.Cells.Item(31, 2).Value = 0.109871145782648
$ws.Cells.Item(32, 1).Value = 0.1899092163002573
$ws.Cells.Item(32, 2).Value = 0.1633960107576634
$ws.Cells.Item(33, 1).Value = 0.1544143112517077
$ws.Cells.Item(33, 2).Value = 0.1521671477337114
$ws.Cells.Item(34, 1).Value = 0.1620024238312819
$ws.Cells.Item(34, 2).Value = 0.09854932069276814
$ws.Cells.Item(35, 1).Value = 0.2298050314572448
$ws.Cells.Item(35, 2).Value = 0.1174619988801536
$ws.Cells.Item(36, 1).Value = 0.124339171062291
$ws.Cells.Item(36, 2).Value = 0.06527680819232659
$ws.Cells.Item(37, 1).Value = 0.1530094231343781
$ws.Cells.Item(37, 2).Value = 0.0594514946830491
$ws.Cells.Item(38, 1).Value = 0.4254958108358347
$ws.Cells.Item(38, 2).Value = 0.2985472568868775
$ws.Cells.Item(39, 1).Value = -0.06040077496932382
$ws.Cells.Item(39, 2).Value = -0.1351409093824841
$ws.Cells.Item(40, 1).Value = 0.2908894226902992
$ws.Cells.Item(40, 2).Value = 0.1493034014773559
$ws.Cells.Item(41, 1).Value = -0.05419303064189934
$ws.Cells.Item(41, 2).Value = -0.0361547641907457
$ws.Cells.Item(42, 1).Value = 0.20347571466625
$ws.Cells.Item(42, 2).Value = 0.1524919071796784
$ws.Cells.Item(43, 1).Value = 0.2815969845777591
$ws.Cells.Item(43, 2).Value = 0.2006684264876506
$ws.Cells.Item(44, 1).Value = -0.1021377025255646
$ws.Cells.Item(44, 2).Value = -0.08552601713209519
$ws.Cells.Item(45, 1).Value = -0.1517115856352539
$ws.Cells.Item(45, 2).Value = -0.1211501334453183
$ws.Cells.Item(46, 1).Value = -0.1933227828001434
$ws.Cells.Item(46, 2).Value = -0.1833007994503395
$ws.Cells.Item(47, 1).Value = -0.1950844104713035
$ws.Cells.Item(47, 2).Value = -0.1861946794937265
$ws.Cells.Item(48, 1).Value = -0.2294127560370921
$ws.Cells.Item(48, 2).Value = -0.2076230289752643
$ws.Cells.Item(49, 1).Value = -0.2249425548125457
$ws.Cells.Item(49, 2).Value = -0.2121478112427399
$ws.Cells.Item(50, 1).Value = -0.1580933802189264
$ws.Cells.Item(50, 2).Value = -0.1531621708789696
$ws.Cells.Item(51, 1).Value = -0.2276570503940634
$ws.Cells.Item(51, 2).Value = -0.2319235299447909
$ws.Cells.Item(52, 1).Value = -0.2276570503940634
$ws.Cells.Item(52, 2).Value = -0.2319235299447909
$ws.Cells.Item(53, 1).Value = -0.2056498981543013
$ws.Cells.Item(53, 2).Value = -0.181458774673453
$ws.Cells.Item(54, 1).Value = -0.2423704189489711
$ws.Cells.Item(54, 2).Value = -0.2220423729941434
$ws.Cells.Item(55, 1).Value = -0.169337150738703
$ws.Cells.Item(55, 2).Value = -0.1623871774418484
$ws.Cells.Item(56, 1).Value = -0.1733844404389098
$ws.Cells.Item(56, 2).Value = -0.1743746655451956
$ws.Cells.Item(57, 1).Value = -0.2285973503574921
$ws.Cells.Item(57, 2).Value = -0.1793497156628741
$ws.Cells.Item(58, 1).Value = -0.2298742512182779
$ws.Cells.Item(58, 2).Value = -0.225510845559846
$ws.Cells.Item(59, 1).Value = -0.2745803355222415
$ws.Cells.Item(59, 2).Value = -0.2543640517107871
$ws.Cells.Item(60, 1).Value = -0.3102167032643084
$ws.Cells.Item(60, 2).Value = -0.287873207325772
$ws.Cells.Item(61, 1).Value = -0.2307152775426941
$ws.Cells.Item(61, 2).Value = -0.1983744637806665
$ws.Cells.Item(62, 1).Value = -0.1889625217605214
$ws.Cells.Item(62, 2).Value = -0.1359259201719744
$ws.Cells.Item(63, 1).Value = -0.3921623931593579
$ws.Cells.Item(63, 2).Value = -0.3994519115559262
$ws.Cells.Item(64, 1).Value = -0.2899266149721538
$ws.Cells.Item(64, 2).Value = -0.2606559396775235
$ws.Cells.Item(65, 1).Value = -0.3621125402306343
$ws.Cells.Item(65, 2).Value = -0.3245313566708664
$ws.Cells.Item(66, 1).Value = -0.1376553587456776
$ws.Cells.Item(66, 2).Value = -0.1374100644147727
$ws.Cells.Item(67, 1).Value = -0.1714767058755454
$ws.Cells.Item(67, 2).Value = -0.1503849287351964
